$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "TC01/TC02" columns (now C),
# shifting the old C/D columns to D/E, to hold Locator / #fname / #lname
# values used by the new "fill form" lookup.
$ws.Columns.Item(3).Insert()

$ws.Range("C1").Value = "Locator"
$ws.Range("C2").Value = "#fname"
$ws.Range("C3").Value = "#lname"

# Match the width already used for column B.
$ws.Columns.Item(3).ColumnWidth = 9.5

# Leave the same selection state Excel would land on after this edit.
$ws.Range("D3").Select()
